# Natmi following Dr Hou advice
# Add "ECs" as a third cluster to the Sending/Target cluster cross-join,
# reordering shared strings (ECs moves before FAPs) and expanding the
# data grid from a 2x2 to a 3x3 sending x target cluster matrix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs", "Vtn", "Plaur", "ECs", 2, 0.6666666666666666, 7.134618, 21.403854, 0.0965317920926077, 0.0965317920926077, 3, 1, 43.91845900000001, 131.755377, 0.8150909120558799, 0.81509091205588, 313.341428113662, 2820.072853022958, 0.07868218645915219, 0.0786821864591522),
    @("ECs", "Vtn", "Plaur", "FAPs", 2, 0.6666666666666666, 7.134618, 21.403854, 0.0965317920926077, 0.0965317920926077, 3, 1, 2.100310333333333, 6.300930999999999, 0.03898005312975703, 0.03898005312975703, 14.984911909786, 134.864207188074, 0.003762814384480508, 0.003762814384480508),
    @("ECs", "Vtn", "Plaur", "sCs", 2, 0.6666666666666666, 7.134618, 21.403854, 0.0965317920926077, 0.0965317920926077, 3, 1, 7.8629, 23.5887, 0.1459290348143631, 0.1459290348143631, 56.0987878722, 504.8890908498, 0.01408679124897501, 0.01408679124897501),
    @("FAPs", "Vtn", "Plaur", "ECs", 3, 1, 17.50798033333334, 52.52394100000001, 0.2368839813846793, 0.2368839813846794, 3, 1, 43.91845900000001, 131.755377, 0.8150909120558799, 0.81509091205588, 768.9235164423065, 6920.311647980759, 0.1930819804382664, 0.1930819804382664),
    @("FAPs", "Vtn", "Plaur", "FAPs", 3, 1, 17.50798033333334, 52.52394100000001, 0.2368839813846793, 0.2368839813846794, 3, 1, 2.100310333333333, 6.300930999999999, 0.03898005312975703, 0.03898005312975703, 36.77219200989678, 330.949728089071, 0.009233750179963175, 0.009233750179963175),
    @("FAPs", "Vtn", "Plaur", "sCs", 3, 1, 17.50798033333334, 52.52394100000001, 0.2368839813846793, 0.2368839813846794, 3, 1, 7.8629, 23.5887, 0.1459290348143631, 0.1459290348143631, 137.6634985629667, 1238.9714870667, 0.03456825076644981, 0.03456825076644982),
    @("sCs", "Vtn", "Plaur", "ECs", 3, 1, 49.26691733333334, 147.800752, 0.6665842265227129, 0.666584226522713, 3, 1, 43.91845900000001, 131.755377, 0.8150909120558799, 0.81509091205588, 2163.72708896039, 19473.54380064351, 0.5433267451584612, 0.5433267451584615),
    @("sCs", "Vtn", "Plaur", "FAPs", 3, 1, 49.26691733333334, 147.800752, 0.6665842265227129, 0.666584226522713, 3, 1, 2.100310333333333, 6.300930999999999, 0.03898005312975703, 0.03898005312975703, 103.4758155666791, 931.282340100112, 0.02598348856531334, 0.02598348856531334),
    @("sCs", "Vtn", "Plaur", "sCs", 3, 1, 49.26691733333334, 147.800752, 0.6665842265227129, 0.666584226522713, 3, 1, 7.8629, 23.5887, 0.1459290348143631, 0.1459290348143631, 387.3808443002667, 3486.4275987024, 0.09727399279893825, 0.09727399279893828)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $row = $rows[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $row.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}

Write-Host "Wrote" $rows.Count "data rows starting at row" $startRow
